$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Planilha1")
if ($ws -eq $null) { $ws = $wb.ActiveSheet }

# Insert a new row at row 26 (everything from the old row 26 downward shifts down by one row),
# carrying formatting down from the row above, same as Excel's default "Insert" behaviour.
$ws.Rows.Item(26).Insert()

# Populate the newly inserted row with the new error-code entry.
# Column A holds the numeric code but stored as TEXT (like every other code in this column,
# e.g. "18", "19", "24"), so force the text number format before assigning the value.
$ws.Range("A26").NumberFormat = "@"
$ws.Range("A26").Value = "25"
$ws.Range("B26").Value = "Dia semana inválido"
$ws.Range("C26").Value = "Tentativa de abastecimento em um dia de semana não permitido"

# Restore the selection to match the saved view state after the edit.
[void]$ws.Range("B27").Select()
